# Update cryptos list worksheet with latest price/volume data
# (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the "General" (unstyled)
# appearance of the cell so numeric-looking strings (e.g. "354.72") are not
# auto-converted into actual numbers by Excel.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "52.050.37"
$ws.Range("E2").Value2 = "  +0.77%  "
Set-TextValue $ws.Range("D3") "3.009.43"
$ws.Range("E3").Value2 = "  +3.21%  "
$ws.Range("E4").Value2 = "  -0.04%  "
Set-TextValue $ws.Range("D5") "354.72"
$ws.Range("E5").Value2 = "  +0.10%  "
Set-TextValue $ws.Range("D6") "107.16"
$ws.Range("E6").Value2 = "  -2.53%  "
$ws.Range("E7").Value2 = "  -1.45%  "
$ws.Range("E8").Value2 = "  +0.19%  "
Set-TextValue $ws.Range("D9") "0.613"
$ws.Range("E9").Value2 = "  -2.46%  "
Set-TextValue $ws.Range("D10") "38.17"
$ws.Range("E10").Value2 = "  -2.65%  "
$ws.Range("E11").Value2 = "  +2.39%  "
Set-TextValue $ws.Range("D12") "0.0858"
$ws.Range("E12").Value2 = "  -3.19%  "
Set-TextValue $ws.Range("D13") "19.05"
$ws.Range("E13").Value2 = "  -3.24%  "
Set-TextValue $ws.Range("D14") "3.476.28"
$ws.Range("E14").Value2 = "  +2.86%  "
$ws.Range("E15").Value2 = "  -3.24%  "
Set-TextValue $ws.Range("D16") "3.020.89"
$ws.Range("E16").Value2 = "  +4.17%  "
$ws.Range("E17").Value2 = "  +4.24%  "
Set-TextValue $ws.Range("D18") "52.122.31"
$ws.Range("E18").Value2 = "  +0.80%  "
$ws.Range("E19").Value2 = "  +3.93%  "
Set-TextValue $ws.Range("D20") "7.47"
$ws.Range("E20").Value2 = "  -1.16%  "
Set-TextValue $ws.Range("D21") "13.59"
$ws.Range("E21").Value2 = "  -2.61%  "
Set-TextValue $ws.Range("D22") "0.0₃0974"
$ws.Range("E22").Value2 = "  -0.61%  "
Set-TextValue $ws.Range("D23") "69.21"
$ws.Range("E23").Value2 = "  -2.23%  "
Set-TextValue $ws.Range("D24") "264.01"
$ws.Range("E24").Value2 = "  -2.00%  "
$ws.Range("E25").Value2 = "  -2.89%  "
$ws.Range("E26").Value2 = "  -1.59%  "
Set-TextValue $ws.Range("D27") "27.03"
$ws.Range("E27").Value2 = "  +0.02%  "
$ws.Range("E28").Value2 = "  +0.05%  "
Set-TextValue $ws.Range("D29") "7.46"
$ws.Range("E29").Value2 = "  +0.96%  "
$ws.Range("E30").Value2 = "  -0.80%  "
Set-TextValue $ws.Range("D31") "6.43"
$ws.Range("E31").Value2 = "  +5.84%  "
$ws.Range("E32").Value2 = "  -3.60%  "
Set-TextValue $ws.Range("D33") "36.14"
$ws.Range("E33").Value2 = "  -4.98%  "
$ws.Range("E34").Value2 = "  +15.41%  "
Set-TextValue $ws.Range("D35") "51.19"
$ws.Range("E35").Value2 = "  -2.19%  "
Set-TextValue $ws.Range("D36") "0.0438"
$ws.Range("E36").Value2 = "  -0.61%  "
Set-TextValue $ws.Range("D37") "0.998"
$ws.Range("E37").Value2 = "  -0.10%  "
Set-TextValue $ws.Range("D38") "3.35"
$ws.Range("E38").Value2 = "  +3.54%  "
$ws.Range("E39").Value2 = "  +4.10%  "
Set-TextValue $ws.Range("D40") "1.97"
$ws.Range("E40").Value2 = "  -2.20%  "
Set-TextValue $ws.Range("D41") "17.66"
$ws.Range("E41").Value2 = "  -4.27%  "
$ws.Range("E42").Value2 = "  -1.47%  "
$ws.Range("B43").Value2 = "EnergySwap"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D43") "23.02"
$ws.Range("E43").Value2 = "  +1.01%  "
$ws.Range("B44").Value2 = "Monero"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D44") "124.39"
$ws.Range("E44").Value2 = "  +2.04%  "
Set-TextValue $ws.Range("D45") "2.13"
$ws.Range("E45").Value2 = "  -1.57%  "
Set-TextValue $ws.Range("D46") "2.124.33"
$ws.Range("E46").Value2 = "  -0.34%  "
Set-TextValue $ws.Range("D47") "3.34"
$ws.Range("E47").Value2 = "  -2.95%  "
$ws.Range("E48").Value2 = "  -6.42%  "
Set-TextValue $ws.Range("D49") "3.303.22"
$ws.Range("E49").Value2 = "  +2.94%  "
Set-TextValue $ws.Range("D50") "0.243"
$ws.Range("E50").Value2 = "  -3.20%  "
Set-TextValue $ws.Range("D51") "0.0335"
$ws.Range("E51").Value2 = "  +0.78%  "
